# Update "想去人数" (interest count) figures in the 广州-漫展信息 workbook.
# The workbook has 4 sheets: 展览, 演出, 本地生活, 全部类型 (the latter being an
# aggregate view of the other three). Only column F values change.

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览" = @{
        "F2"  = 43
        "F5"  = 987
        "F6"  = 373
        "F7"  = 605
        "F10" = 133
        "F11" = 1360
        "F12" = 3017
        "F13" = 454
        "F14" = 1647
        "F18" = 1400
        "F21" = 1130
        "F22" = 3
        "F24" = 10
        "F25" = 3513
        "F26" = 697
        "F27" = 561
        "F28" = 1554
    }
    "演出" = @{
        "F5" = 32
        "F6" = 12
        "F7" = 54
        "F8" = 25
    }
    "本地生活" = @{
        "F3" = 11
    }
    "全部类型" = @{
        "F2"  = 43
        "F5"  = 11
        "F9"  = 32
        "F10" = 12
        "F11" = 54
        "F13" = 25
        "F15" = 987
        "F16" = 373
        "F17" = 605
        "F20" = 133
        "F21" = 1360
        "F22" = 3017
        "F23" = 454
        "F24" = 1647
        "F28" = 1400
        "F33" = 1130
        "F34" = 3
        "F36" = 10
        "F37" = 3513
        "F38" = 697
        "F39" = 561
        "F40" = 1554
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellUpdates = $updates[$sheetName]
    foreach ($addr in $cellUpdates.Keys) {
        $ws.Range($addr).Value = $cellUpdates[$addr]
    }
}
